$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 480
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 480
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 480
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -830
# row 98
$ws.Range("H98").Value = 555.3125
$ws.Range("I98").Value = 563.3570999999999
$ws.Range("K98").Value = 563.3570999999999
$ws.Range("M98").Value = 934.6429000000001
# row 100
$ws.Range("H100").Value = 1886.9166
$ws.Range("I100").Value = 1972.091
$ws.Range("K100").Value = 1972.091
$ws.Range("M100").Value = -1431.091
# row 113
$ws.Range("H113").Value = 2831.6667
$ws.Range("I113").Value = 2831.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2831.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 422.3332999999998
$ws.Range("N113").ClearContents()
# row 122
$ws.Range("H122").Value = 555.3125
$ws.Range("I122").Value = 563.3570999999999
$ws.Range("K122").Value = 1690.0713
$ws.Range("M122").Value = 759.9287000000002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 6402.6
$ws.Range("I61").Value = 7666.3335
$ws.Range("K61").Value = 7666.3335
$ws.Range("M61").Value = -7454.3335
# row 122
$ws.Range("H122").Value = 2002.3334
$ws.Range("I122").Value = 1503.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4510.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2060.5
$ws.Range("N122").Value = -13900
# row 132
$ws.Range("H132").Value = 1101
$ws.Range("I132").Value = 1101
$ws.Range("K132").Value = 3303
$ws.Range("M132").Value = -773
# row 136
$ws.Range("H136").Value = 6402.6
$ws.Range("I136").Value = 7666.3335
$ws.Range("K136").Value = 22999.0005
$ws.Range("M136").Value = -20449.0005
# row 139
$ws.Range("H139").Value = 69999.5
$ws.Range("J139").Value = 69999.5
$ws.Range("L139").Value = 69999.5
$ws.Range("N139").Value = -80279.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 1889.5
$ws.Range("I134").Value = 1889.5
$ws.Range("K134").Value = 5668.5
$ws.Range("M134").Value = -3133.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 1027.7142
$ws.Range("I22").Value = 533
$ws.Range("J22").Value = 1398.75
$ws.Range("K22").Value = 533
$ws.Range("L22").Value = 1398.75
$ws.Range("M22").Value = -183
$ws.Range("N22").Value = -2098.75
# row 31
$ws.Range("H31").Value = 2915.5557
$ws.Range("I31").Value = 2915.5557
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2915.5557
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2620.5557
$ws.Range("N31").ClearContents()
# row 34
$ws.Range("H34").Value = 2915.5557
$ws.Range("I34").Value = 2915.5557
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2915.5557
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2713.5557
$ws.Range("N34").ClearContents()
# row 62
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# row 65
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# row 86
$ws.Range("H86").Value = 5811455
$ws.Range("I86").Value = 7747217.5
$ws.Range("J86").Value = 4166.6665
$ws.Range("K86").Value = 7747217.5
$ws.Range("L86").Value = 4166.6665
$ws.Range("M86").Value = -7746094.5
$ws.Range("N86").Value = -6412.6665
# row 89
$ws.Range("H89").Value = 5811455
$ws.Range("I89").Value = 7747217.5
$ws.Range("J89").Value = 4166.6665
$ws.Range("K89").Value = 38736087.5
$ws.Range("L89").Value = 20833.3325
$ws.Range("M89").Value = -38730471.5
$ws.Range("N89").Value = -32065.3325
# row 138
$ws.Range("H138").Value = 1427.7142
$ws.Range("J138").Value = 1000
$ws.Range("L138").Value = 1000
$ws.Range("N138").Value = -11280

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("I5").Value = 301
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 903
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -791
$ws.Range("N5").ClearContents()
# row 44
$ws.Range("H44").Value = 691.1429000000001
$ws.Range("I44").Value = 691.1429000000001
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 2073.4287
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -1675.4287
$ws.Range("N44").ClearContents()
# row 55
$ws.Range("H55").Value = 1429.5714
$ws.Range("I55").Value = 402.66666
$ws.Range("J55").Value = 2199.75
$ws.Range("K55").Value = 1207.99998
$ws.Range("L55").Value = 6599.25
$ws.Range("M55").Value = -1030.99998
$ws.Range("N55").Value = -6953.25
# row 118
$ws.Range("H118").Value = 839.25
$ws.Range("I118").Value = 839.25
$ws.Range("K118").Value = 2517.75
$ws.Range("M118").Value = -1274.75
# row 122
$ws.Range("H122").Value = 1118.5
$ws.Range("I122").Value = 624.6667
$ws.Range("J122").Value = 1488.875
$ws.Range("K122").Value = 5622.0003
$ws.Range("L122").Value = 13399.875
$ws.Range("M122").Value = -3172.0003
$ws.Range("N122").Value = -18299.875
# row 135
$ws.Range("I135").Value = 301
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2709
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -174
$ws.Range("N135").ClearContents()
# row 137
$ws.Range("H137").Value = 2287.8572
$ws.Range("I137").Value = 1700
$ws.Range("J137").Value = 2385.8333
$ws.Range("K137").Value = 5100
$ws.Range("L137").Value = 7157.499899999999
$ws.Range("N137").Value = -17357.4999
$ws.Range("M137").Value = 0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 57
$ws.Range("H57").Value = 2950
$ws.Range("I57").Value = 2950
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 2950
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2130
$ws.Range("N57").ClearContents()
# row 122
$ws.Range("H122").Value = 5600
$ws.Range("I122").Value = 5600
$ws.Range("K122").Value = 16800
$ws.Range("M122").Value = -14350
# row 132
$ws.Range("H132").Value = 3021.6667
$ws.Range("I132").Value = 2836.2
$ws.Range("K132").Value = 8508.599999999999
$ws.Range("M132").Value = -5978.599999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 842.7857
$ws.Range("I55").Value = 339.4
$ws.Range("J55").Value = 1122.4445
$ws.Range("K55").Value = 339.4
$ws.Range("L55").Value = 1122.4445
$ws.Range("M55").Value = -166.4
$ws.Range("N55").Value = -1468.4445
# row 122
$ws.Range("H122").Value = 5985
$ws.Range("I122").Value = 5985
$ws.Range("K122").Value = 17955
$ws.Range("M122").Value = -15505

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 105
$ws.Range("H105").Value = 36250
$ws.Range("J105").Value = 36250
$ws.Range("L105").Value = 36250
$ws.Range("N105").Value = -43238
# row 126
$ws.Range("H126").Value = 1862.5
$ws.Range("J126").Value = 1725
$ws.Range("L126").Value = 5175
$ws.Range("N126").Value = -10115
# row 132
$ws.Range("H132").Value = 3340.3635
$ws.Range("J132").Value = 3112.5
$ws.Range("L132").Value = 9337.5
$ws.Range("N132").Value = -14397.5
# row 136
$ws.Range("H136").Value = 3810
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

Write-Host "Edits applied"